$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = 25.8400000000006
$ws.Range("G2").Value = 0.000128378529096751
$ws.Range("H2").Value = 0.003258625825198164
$ws.Range("K2").Value = 4.36576024717086
$ws.Range("L2").Value = "[1.9301836425729402, 6.801336851768779]"
$ws.Range("M2").Value = 0.0004692003707209569
$ws.Range("N2").Value = 0.0009384007414419138
$ws.Range("O2").Value = -1.42142130072331
$ws.Range("P2").Value = "[-2.1132635267390807, -0.729579074707539]"
$ws.Range("Q2").Value = [double]"6.304169505977342e-05"
$ws.Range("R2").Value = 0.0001260833901195468
$ws.Range("S2").Value = 13.64314009029865
$ws.Range("T2").Value = "[12.220562583149047, 15.065717597448252]"
$ws.Range("W2").Value = 5.845685685685822
$ws.Range("X2").Value = 3.000440440440509
$ws.Range("Y2").Value = 8.690930930931136

# --- Row 3 updates ---
$ws.Range("B3").Value = 0
$ws.Range("G3").Value = 0.01088044110989639
$ws.Range("H3").Value = 0.06643600926369818
$ws.Range("I3").Value = [double]"2.315581509870768e-07"
$ws.Range("K3").Value = 5.222625362313411
$ws.Range("L3").Value = "[0.9287875602336388, 9.516463164393183]"
$ws.Range("M3").Value = 0.01739891774252045
$ws.Range("N3").Value = 0.01739891774252045
$ws.Range("O3").Value = 0.5094474573388847
$ws.Range("P3").Value = "[-0.5534737808126167, 1.5723686954903862]"
$ws.Range("Q3").Value = 0.3456388677100648
$ws.Range("R3").Value = 0.3456388677100648
$ws.Range("S3").Value = 13.33983611157445
$ws.Range("T3").Value = "[10.95409299237457, 15.725579230774336]"
$ws.Range("W3").Value = 20.21621621621622
$ws.Range("X3").Value = 16.4944944944945
$ws.Range("Y3").Value = 23.93793793793794
